$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Shift the existing "Quartiers" values (column C, rows 3:15) down by one row,
# working from the bottom up so no value is overwritten before it is copied,
# then insert the new entry at the top of the list.
for ($r = 15; $r -ge 3; $r--) {
    $ws.Cells.Item($r + 1, 3).Value2 = $ws.Cells.Item($r, 3).Value2
}
$ws.Cells.Item(3, 3).Value2 = "AffaireDroit"

# Update the active selection to match the post-edit state.
$ws.Range("E22").Select()
